# codeforIATI/codelists ReportingOrganisationGroup.xlsx regeneration:
# the "codeforiati:group-code" and "codeforiati:group-name" columns swap
# order (code now precedes name, matching the code/name convention used
# by the other columns). Every row's data stays the same - only the two
# right-most columns (D and E) trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow  = $used.Row
$lastRow   = $firstRow + $used.Rows.Count - 1
$lastCol   = $used.Column + $used.Columns.Count - 1
$groupNameCol = $lastCol - 1   # "codeforiati:group-name" column (D)
$groupCodeCol = $lastCol       # "codeforiati:group-code" column (E)

$rng = $ws.Range($ws.Cells.Item($firstRow, $groupNameCol), $ws.Cells.Item($lastRow, $groupCodeCol))
$vals = $rng.Value()

$rowCount = $vals.GetLength(0)
$swapped = New-Object 'object[,]' $rowCount,2
for ($i = 1; $i -le $rowCount; $i++) {
    $swapped[$i - 1, 0] = $vals[$i, 2]
    $swapped[$i - 1, 1] = $vals[$i, 1]
}

$rng.Value = $swapped
